# MegaSena_edt.xlsx - "Limpar Tudo" purple-residue fix
#
# Behaviour being reproduced:
#  1. The previous "last draws" highlight (purple clover style, s="2" on col A
#     and s="3" on cols B:G) that was sitting on rows 381:390 is cleared back
#     to plain/no style - this is the actual bug fix ("Corrige limpeza de
#     trevos no botao Limpar Tudo").
#  2. Four new Mega-Sena draws (contests 2934-2937) are appended as rows
#     391:394, re-using that same highlight style for the newest block.
#  3. The sheet dimension / selection follow the newly used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the stale purple-residue formatting from the old "latest draws" block ---
$ws.Range("A381:G390").ClearFormats()

# --- 2. Append the new draws (contests 2934-2937) ---
$newDraws = @(
    @(2934, 9, 17, 23, 26, 33, 59),
    @(2935, 9, 18, 28, 34, 38, 57),
    @(2936, 4, 7, 9, 15, 29, 32),
    @(2937, 12, 17, 26, 34, 44, 52)
)

$startRow = 391
for ($i = 0; $i -lt $newDraws.Count; $i++) {
    $r = $startRow + $i
    $row = $newDraws[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# --- 3. Update the active view to match where the latest draws now live ---
$ws.Activate()
$ws.Range("B391:G394").Select()
